# "update nav interne review"
# Fills in the "verwerking" (column G) feedback/response text for several
# rows on the "Blad1" review sheet, updates one existing response's wording,
# and moves the active selection to G20 (where the reviewer was last working).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")
$ws.Activate()

# Reword the existing response in row 11 (already filled in).
$ws.Range("G11").Value = "verwerkt in uml en daarmee ook in modeldocument (na generatie objectcat)"

# Row 12: previously empty G cell -> multi-line response.
$ws.Range("G12").Value = "verwerkt in uml`nverwerkt in xsd`nverwerkt in changelog"

# Rows 13, 14 and 17 all get the same short response.
$ws.Range("G13").Value = "is verwerkt in waardelijst excel"
$ws.Range("G14").Value = "is verwerkt in waardelijst excel"
$ws.Range("G17").Value = "is verwerkt in waardelijst excel"

# Row 15: new cell, needs wrap-text styling (matches column D's style) and a
# longer, multi-line response with a trailing blank line.
$ws.Range("G15").WrapText = $true
$ws.Range("G15").Value = "changlog item 39 is verwijderd.`nVerwerkt in UML: Constraints zijn aangepast/herformuleerd op basis van voorstel.`nVerwerkt in changelog: herformulering van constraints op GebiedsinformatieAanvraag`n"

# Row 18 and 19: previously empty G cells -> short responses.
$ws.Range("G18").Value = "is verwerkt in UML diagrammen en verwijderd uit model."
$ws.Range("G19").Value = "is verwerkt in extraRegels excel"

# Row 20: new cell.
$ws.Range("G20").Value = "verwerkt. Dit changelog item is nu item 42 geworden"

# Move/restore the selection to G20, matching where the reviewer left off.
$ws.Range("G20").Select() | Out-Null
